# Apply "Updates to Variable Naming Convention" edit:
# Cell J4 (Attachments column, row for JR12356) text is corrected from
# "...OMES Background Check Release HPABCD.pdf" to "...OMES Background Check Release HP.pdf"
# Also reselect cell J5 to match the saved cursor position in the target file.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OfferLetter")

$ws.Range("J4").Value = "C:\Users\55649C\Documents\Data\OMES Background Check Release HP.pdf"

$ws.Range("J5").Select()

$wb.Save()
